# Fruta / hortaliza, semanal
# Adds the latest weekly price-report block (fecha 44706) for
# "Terminal La Palmera de La Serena - Piña" at the top of the
# date-descending data range, pushing the existing rows down by
# one 4-row block (Especial/Primera/Segunda/Tercera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 4 blank rows right above the first data row (907).
# Excel inherits formatting (e.g. the date style on column D) from the
# row above automatically, matching the rest of the column.
$ws.Range("A907:A910").EntireRow.Insert()

$newRows = @(
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44706, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Especial", 216, 18000, 18500, 18250, "`$/caja 10 unidades", "Ecuador", 1825, 10),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44706, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Primera", 216, 18000, 18500, 18250, "`$/caja 12 unidades", "Ecuador", 1521, 12),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44706, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Segunda", 216, 18000, 18500, 18250, "`$/caja 14 unidades", "Ecuador", 1304, 14),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44706, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Tercera", 216, 18000, 18500, 18250, "`$/caja 16 unidades", "Ecuador", 1141, 16)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = 907 + $i
    $vals = $newRows[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
}
